# ---------------------------------------------------------------------------
# Helper: convert a "RRGGBB" hex string into the BGR-packed integer that the
# PowerPoint object model expects for ColorFormat/ThemeColor .RGB assignments
# (VBA's RGB() encodes as 0x00BBGGRR).
# ---------------------------------------------------------------------------
function Convert-HexToBGR {
    param([string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16's table picked up the built-in "Medium Style 2 - Accent 1"
#    table style instead of the deck's custom style - re-apply the correct
#    built-in style GUID.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B2E79C4F-6E93-4170-8624-EC24689A74E3}")
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme colour scheme ("Integral") was swapped out for the
#    default Office theme palette - push the Office theme colours onto the
#    presentation's theme colour scheme.
# ---------------------------------------------------------------------------
$officeThemeColors = [ordered]@{
    1  = "000000"   # dk1
    2  = "FFFFFF"   # lt1
    3  = "44546A"   # dk2
    4  = "E7E6E6"   # lt2
    5  = "5B9BD5"   # accent1
    6  = "ED7D31"   # accent2
    7  = "A5A5A5"   # accent3
    8  = "FFC000"   # accent4
    9  = "4472C4"   # accent5
    10 = "70AD47"   # accent6
    11 = "0563C1"   # hlink
    12 = "954F72"   # folHlink
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $officeThemeColors.Keys) {
    $tcs.Colors($idx).RGB = Convert-HexToBGR $officeThemeColors[$idx]
}
